$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RCommands")

$ws.Range("A119").Value = "General"
$ws.Range("B119").Value = "nrow()"
$ws.Range("C119").Value = "nrow(Data frame A)"
$ws.Range("D119").Value = "get the rows of A"

$ws.Range("A120").Value = "General"
$ws.Range("B120").Value = "ncol()"
$ws.Range("C120").Value = "ncol(Data frame A)"
$ws.Range("D120").Value = "get the columns of A"

$ws.Range("A121").Value = "Package: reshape2"
$ws.Range("B121").Value = "melt()"
$ws.Range("C121").Value = 'mtcars$car <- rownames(mtcars); mtcarsMelt <- melt(mtcars)'
$ws.Range("D121").Value = "reshape data frame mtcars into 3 columns: [id, variables, value]"

$ws.Activate()
$ws.Range("C26:C30").Select()

Write-Output "done"
